$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A
$lastRow = $ws.Cells.Item(1,1).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    $d = [DateTime]::FromOADate($serial)
    $d2 = $d.AddMonths(1)
    $d3 = Get-Date -Year $d2.Year -Month $d2.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $newSerial = $d3.ToOADate()
    $cell.Value = $newSerial
}
